$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "isMgr" column (column D) entirely
$ws.Columns.Item(4).Delete() | Out-Null

# Repurpose the "Designation" column (column B) into a "Location" column
$ws.Range("B1").Value = "Location"
$ws.Range("B2").Value = "Hi TechCity"
$ws.Range("B3").Value = "Madhabpur"

# Match the selection shown in the saved workbook
$ws.Range("B3").Select() | Out-Null

# Rename the existing sheet and add a new, empty "Sheet2" after it
$ws.Name = "empDetails"
$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "Sheet2"

# Keep the original sheet as the active/selected tab
$ws.Activate() | Out-Null
